# Suppl. Table 5 update: row 10 ("Secondary" education) values change
# because the table no longer uses age-standardized estimates for the
# age-group breakdown. Values are stored as text in this sheet (several
# other cells hold non-numeric placeholders like "-"), so force the
# range to Text format before assigning the new figures to keep the
# cells as text rather than letting Excel auto-convert them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row10 = $ws.Range("B10:L10")
$row10.NumberFormat = "@"

$ws.Range("B10").Value = "0.29"
$ws.Range("C10").Value = "0.5"
$ws.Range("D10").Value = "0.4"
$ws.Range("E10").Value = "0.49"
$ws.Range("F10").Value = "0.41"
$ws.Range("G10").Value = "0.55"
$ws.Range("H10").Value = "0.44"
$ws.Range("I10").Value = "0.62"
$ws.Range("J10").Value = "0.66"
$ws.Range("K10").Value = "0.58"
$ws.Range("L10").Value = "0.58"
